$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "bank" column (B) data rows (2-156) currently hold the Thai value
# "ธนาคาร" ("bank"). Update every data row to "ธนาคาร UOB" while leaving
# the header in B1 ("bank") untouched.
$dataRange = $ws.Range("B2:B156")
$dataRange.Value = "ธนาคาร UOB"

# Reflect the new selection left behind by the edit (B2:B156, active cell B2).
$dataRange.Select()
